$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '59.213.33'
$ws.Range('E2').Value = '  -0.33%  '

# Row 3
$ws.Range('D3').Value = '2.635.37'
$ws.Range('E3').Value = '  -0.75%  '

# Row 4
$ws.Range('E4').Value = '  -0.11%  '

# Row 5
$ws.Range("D5").Value = "'527.66"
$ws.Range('E5').Value = '  +1.38%  '

# Row 6
$ws.Range("D6").Value = "'144.58"
$ws.Range('E6').Value = '  -1.70%  '

# Row 7
$ws.Range('E7').Value = '  -0.19%  '

# Row 8
$ws.Range("D8").Value = "'0.570"
$ws.Range('E8').Value = '  -0.05%  '

# Row 9
$ws.Range("D9").Value = "'6.66"
$ws.Range('E9').Value = '  -3.76%  '

# Row 10
$ws.Range('E10').Value = '  +1.58%  '

# Row 11
$ws.Range('E11').Value = '  +0.74%  '

# Row 12
$ws.Range('E12').Value = '  +0.99%  '

# Row 13
$ws.Range('D13').Value = '3.101.02'
$ws.Range('E13').Value = '  -0.65%  '

# Row 14
$ws.Range('D14').Value = '59.168.24'
$ws.Range('E14').Value = '  -0.35%  '

# Row 15
$ws.Range("D15").Value = "'20.99"
$ws.Range('E15').Value = '  -0.23%  '

# Row 16
$ws.Range('B16').Value = 'WrappedEther'
$ws.Range('C16').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D16').Value = '2.679.36'
$ws.Range('E16').Value = '  +1.23%  '

# Row 17
$ws.Range('B17').Value = 'ShibaInu'
$ws.Range('C17').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D17").Value = "'0.0000137"
$ws.Range('E17').Value = '  +0.66%  '

# Row 18
$ws.Range("D18").Value = "'340.92"
$ws.Range('E18').Value = '  +0.45%  '

# Row 19
$ws.Range('E19').Value = '  +0.51%  '

# Row 20
$ws.Range('E20').Value = '  +2.57%  '

# Row 21
$ws.Range("D21").Value = "'6.34"
$ws.Range('E21').Value = '  +0.74%  '

# Row 22
$ws.Range("D22").Value = "'0.999"
$ws.Range('E22').Value = '  +0.29%  '

# Row 23
$ws.Range("D23").Value = "'65.46"
$ws.Range('E23').Value = '  +3.68%  '

# Row 24
$ws.Range('E24').Value = '  +1.56%  '

# Row 25
$ws.Range("D25").Value = "'0.168"
$ws.Range('E25').Value = '  +0.14%  '

# Row 26
$ws.Range("D26").Value = "'0.998"
$ws.Range('E26').Value = '  -0.51%  '

# Row 27
$ws.Range("D27").Value = "'7.24"
$ws.Range('E27').Value = '  +1.27%  '

# Row 28
$ws.Range('E28').Value = '  -0.74%  '

# Row 29
$ws.Range('B29').Value = 'Aptos'
$ws.Range('C29').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D29").Value = "'6.46"
$ws.Range('E29').Value = '  -3.12%  '

# Row 30
$ws.Range('B30').Value = 'USDe'
$ws.Range('C30').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D30").Value = "'0.998"
$ws.Range('E30').Value = '  -0.06%  '

# Row 31
$ws.Range('E31').Value = '  +2.16%  '

# Row 32
$ws.Range('E32').Value = '  +0.95%  '

# Row 33
$ws.Range("D33").Value = "'150.06"
$ws.Range('E33').Value = '  +0.40%  '

# Row 34
$ws.Range("D34").Value = "'4.19"
$ws.Range('E34').Value = '  +0.93%  '

# Row 35
$ws.Range('E35').Value = '  +0.32%  '

# Row 36
$ws.Range("D36").Value = "'0.904"

# Row 37
$ws.Range("D37").Value = "'0.868"
$ws.Range('E37').Value = '  -1.42%  '

# Row 38
$ws.Range('E38').Value = '  +0.18%  '

# Row 39
$ws.Range("D39").Value = "'36.58"
$ws.Range('E39').Value = '  -0.86%  '

# Row 40
$ws.Range('E40').Value = '  +1.97%  '

# Row 41
$ws.Range('E41').Value = '  -0.15%  '

# Row 42
$ws.Range('E42').Value = '  +0.24%  '

# Row 43
$ws.Range('B43').Value = 'Mantle'
$ws.Range('C43').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D43").Value = "'0.601"
$ws.Range('E43').Value = '  -4.34%  '

# Row 44
$ws.Range('B44').Value = 'Bittensor'
$ws.Range('C44').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D44").Value = "'271.93"
$ws.Range('E44').Value = '  -1.55%  '

# Row 45
$ws.Range("D45").Value = "'19.38"
$ws.Range('E45').Value = '  -2.08%  '

# Row 46
$ws.Range('B46').Value = 'WhiteBITCoin'
$ws.Range('C46').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range("D46").Value = "'10.66"
$ws.Range('E46').Value = '  +1.42%  '

# Row 47
$ws.Range('B47').Value = 'Hedera'
$ws.Range('C47').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D47").Value = "'0.0537"
$ws.Range('E47').Value = '  +0.20%  '

# Row 48
$ws.Range('D48').Value = '2.039.23'
$ws.Range('E48').Value = '  -1.06%  '

# Row 49
$ws.Range('B49').Value = 'RenderToken'
$ws.Range('C49').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D49").Value = "'4.78"
$ws.Range('E49').Value = '  -0.47%  '

# Row 50
$ws.Range('B50').Value = 'VeChain'
$ws.Range('C50').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D50").Value = "'0.0230"
$ws.Range('E50').Value = '  +0.10%  '

# Row 51
$ws.Range('E51').Value = '  -0.39%  '
